# datacamp.xlsx — add a new "Understanding Machine Learning" row just above
# the "Machine Learning with scikit-learn" block, pushing the existing rows
# 44-48 down by one (to 45-49), then restore the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 44 (rows 44-48 shift down to 45-49); the new row
# inherits formatting from the row above it, same as Excel's UI insert.
$ws.Rows.Item(44).Insert()

# Populate the new row with the course name + rating, matching the style
# already used by the other "Understanding ..." summary rows (A: s="4",
# I: s="103").
$ws.Range("A44").Value = "Understanding Machine Learning"
$ws.Range("I44").Value = 5

# Leave the selection where the author left it after the edit.
[void]$ws.Range("I45").Select()
